$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and 1h volume change (E) columns.
# D-column updates use a Text number format while writing so that
# numeric-looking strings (e.g. "1.002") are stored as text, matching
# the source data (inline strings), then ClearFormats() restores the
# original (unstyled) cell formatting.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "29.998.57"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +3.61%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.892.48"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +3.59%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.ClearFormats()
$ws.Range("E4").Value = "  +0.37%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "247.74"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +0.23%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +0.26%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4935"
$cell.ClearFormats()
$ws.Range("E7").Value = "  -0.16%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "44.74"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +2.20%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.2934"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +5.25%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.06606"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +3.06%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.895.17"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +3.82%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "16.93"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +0.94%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.07219"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +2.22%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.6738"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +4.39%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "85.89"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +1.97%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "4.837"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +3.35%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "29.999.77"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +3.54%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000007889"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +7.72%  "
$ws.Range("E19").Value = "  +0.34%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "12.83"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +4.72%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "2.142.83"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +4.99%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +0.37%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "4.758"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +4.06%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "5.641"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +5.06%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "9.170"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +3.48%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "147.03"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +1.74%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "130.57"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +0.84%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "16.70"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +1.62%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.955"
$cell.ClearFormats()
$ws.Range("E29").Value = "  +3.48%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.358"
$cell.ClearFormats()
$ws.Range("E30").Value = "  -3.86%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "4.198"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +1.52%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.08697"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +3.97%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.923"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +3.42%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.05076"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +2.94%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.116"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +1.40%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.6986"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +3.67%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.678"
$cell.ClearFormats()
$ws.Range("E37").Value = "  -0.54%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.770"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +2.08%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.209"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -3.97%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.9491"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +0.16%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.01652"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +4.09%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "5.961"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -3.10%  "
$ws.Range("E43").Value = "  +0.08%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.4202"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +2.85%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "102.62"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +1.99%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "7.432"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +3.21%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.1256"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +2.59%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.05737"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +3.87%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "32.66"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +3.08%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "8.233"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +0.86%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.3714"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +2.85%  "
